$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update revised AgTests (F) and AgPosit (G) values for rows 334-377
$ws.Cells.Item(334, 6).Value = 195457
$ws.Cells.Item(334, 7).Value = 3445
$ws.Cells.Item(335, 6).Value = 131254
$ws.Cells.Item(335, 7).Value = 3002
$ws.Cells.Item(336, 6).Value = 101659
$ws.Cells.Item(336, 7).Value = 3366
$ws.Cells.Item(337, 6).Value = 104432
$ws.Cells.Item(337, 7).Value = 2988
$ws.Cells.Item(338, 6).Value = 226806
$ws.Cells.Item(338, 7).Value = 3187
$ws.Cells.Item(339, 6).Value = 657208
$ws.Cells.Item(339, 7).Value = 5480
$ws.Cells.Item(340, 6).Value = 383782
$ws.Cells.Item(340, 7).Value = 3298
$ws.Cells.Item(341, 6).Value = 291522
$ws.Cells.Item(341, 7).Value = 3665
$ws.Cells.Item(342, 6).Value = 179575
$ws.Cells.Item(342, 7).Value = 3073
$ws.Cells.Item(343, 6).Value = 132484
$ws.Cells.Item(343, 7).Value = 2965
$ws.Cells.Item(344, 6).Value = 135698
$ws.Cells.Item(344, 7).Value = 2488
$ws.Cells.Item(345, 6).Value = 291779
$ws.Cells.Item(345, 7).Value = 3315
$ws.Cells.Item(346, 6).Value = 669733
$ws.Cells.Item(346, 7).Value = 4787
$ws.Cells.Item(347, 6).Value = 342166
$ws.Cells.Item(347, 7).Value = 2905
$ws.Cells.Item(348, 6).Value = 232255
$ws.Cells.Item(348, 7).Value = 3250
$ws.Cells.Item(349, 6).Value = 159657
$ws.Cells.Item(349, 7).Value = 2758
$ws.Cells.Item(350, 6).Value = 127245
$ws.Cells.Item(350, 7).Value = 2785
$ws.Cells.Item(351, 6).Value = 150637
$ws.Cells.Item(351, 7).Value = 2827
$ws.Cells.Item(352, 6).Value = 306392
$ws.Cells.Item(352, 7).Value = 3530
$ws.Cells.Item(353, 6).Value = 720596
$ws.Cells.Item(353, 7).Value = 5261
$ws.Cells.Item(354, 6).Value = 308851
$ws.Cells.Item(354, 7).Value = 2831
$ws.Cells.Item(355, 6).Value = 222891
$ws.Cells.Item(355, 7).Value = 3465
$ws.Cells.Item(356, 6).Value = 159978
$ws.Cells.Item(356, 7).Value = 2878
$ws.Cells.Item(357, 6).Value = 138575
$ws.Cells.Item(357, 7).Value = 3032
$ws.Cells.Item(358, 6).Value = 157318
$ws.Cells.Item(358, 7).Value = 2599
$ws.Cells.Item(359, 6).Value = 320363
$ws.Cells.Item(359, 7).Value = 3335
$ws.Cells.Item(360, 6).Value = 744097
$ws.Cells.Item(360, 7).Value = 5105
$ws.Cells.Item(361, 6).Value = 331368
$ws.Cells.Item(361, 7).Value = 2617
$ws.Cells.Item(362, 6).Value = 227642
$ws.Cells.Item(362, 7).Value = 3159
$ws.Cells.Item(363, 6).Value = 187568
$ws.Cells.Item(363, 7).Value = 2748
$ws.Cells.Item(364, 6).Value = 166967
$ws.Cells.Item(364, 7).Value = 2449
$ws.Cells.Item(365, 6).Value = 179636
$ws.Cells.Item(365, 7).Value = 2346
$ws.Cells.Item(366, 6).Value = 336323
$ws.Cells.Item(366, 7).Value = 2822
$ws.Cells.Item(367, 6).Value = 753824
$ws.Cells.Item(367, 7).Value = 3858
$ws.Cells.Item(368, 6).Value = 342922
$ws.Cells.Item(368, 7).Value = 2261
$ws.Cells.Item(369, 6).Value = 232044
$ws.Cells.Item(369, 7).Value = 2562
$ws.Cells.Item(370, 6).Value = 180724
$ws.Cells.Item(370, 7).Value = 2017
$ws.Cells.Item(371, 6).Value = 157639
$ws.Cells.Item(371, 7).Value = 1931
$ws.Cells.Item(372, 6).Value = 175050
$ws.Cells.Item(372, 7).Value = 1813
$ws.Cells.Item(373, 6).Value = 339789
$ws.Cells.Item(373, 7).Value = 2306
$ws.Cells.Item(374, 6).Value = 740998
$ws.Cells.Item(374, 7).Value = 3301
$ws.Cells.Item(375, 6).Value = 337419
$ws.Cells.Item(375, 7).Value = 1795
$ws.Cells.Item(376, 6).Value = 215095
$ws.Cells.Item(376, 7).Value = 2154
$ws.Cells.Item(377, 6).Value = 166869
$ws.Cells.Item(377, 7).Value = 1745

# Append new row 378 for 2021-03-17 data
$ws.Cells.Item(378, 1).Value = 44272
$ws.Cells.Item(378, 2).Value = 344470
$ws.Cells.Item(378, 3).Value = 11088
$ws.Cells.Item(378, 4).Value = 2040
$ws.Cells.Item(378, 5).Value = 8814
$ws.Cells.Item(378, 6).Value = 129942
$ws.Cells.Item(378, 7).Value = 1275

Write-Host "Update complete"
